# Opp to Eng validation and expense req changes
#
# On the "TitleRateSheet" worksheet, the rate-sheet name in row 2 changes
# from "Schedule A" to "DRC - Original". The active selection on that sheet
# also moves from D21 to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TitleRateSheet")

$ws.Range("A2").Value = "DRC - Original"

$ws.Range("D7").Select()
